$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 2. Update column A values (Gen counts -> MaxFES fractions)
$colA = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $colA[$i]
}

# 3. Remove the "Run 50" column (AZ) entirely - "Mean" (was BA) shifts left into AZ
$ws.Range("AZ1:AZ14").EntireColumn.Delete()

# 4. Recompute/update the Mean column (now AZ) values, excluding the removed Run 50 data
$meanVals = @(
    626.76949671,
    626.76949671,
    578.6768499999999,
    379.46039652,
    315.56897951,
    280.96375073,
    264.06963418,
    249.47428349,
    232.62752581,
    220.37986661,
    212.48152546,
    205.60969929,
    199.32369083
)
for ($i = 0; $i -lt $meanVals.Length; $i++) {
    $ws.Cells.Item($i + 2, 52).Value = $meanVals[$i]
}
